$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "C3" = 37.793897434793102
    "D3" = 3.83017279663627
    "E3" = 8.4707680603003208
    "F3" = 0.19240763751553
    "G3" = 21595.055055266399
    "I3" = 4840.1121619750002
    "K3" = 10.814415339865
    "L3" = 2.2716142341359999
    "C4" = 38.401811811407804
    "D4" = 0.042792340346119999
    "E4" = 7.4656456644002702
    "F4" = -0.13636111348630001
    "H4" = 6.3888964136758197
    "J4" = -20.358714243502
    "K4" = 0.11421685021855001
    "L4" = -1.9143110861257
    "C5" = 24.1657376925404
    "D5" = 11.416157692100899
    "E5" = 3.6744753514451798
    "F5" = 0.11188403866740999
    "C6" = 52.958962659692503
    "D6" = 4.5618633647694304
    "C7" = 12.7377964391503
    "D7" = 0.86310045315282002
    "E7" = 1.88461496451433
    "C8" = 327.02641462365801
    "D8" = 0.88994900686869005
    "E8" = 0.13192105195375001
    "F8" = 0.076935711406590002
    "I8" = 18.676063325092802
    "L8" = 19.015173373142598
    "C9" = 66.584132003183399
    "D9" = 0.51438613566207003
    "E9" = 10.0747554392241
    "F9" = 0.02191130172199
    "G9" = 7046.5986898969004
    "I9" = 1066.21136813309
    "K9" = 0.78473809184518994
    "L9" = 0.21420023291314
    "C10" = 39.364083469202498
    "D10" = 1.1038669713736999
    "E10" = 74.207203693571699
    "F10" = 1.13780929159358
    "H10" = 3733.2008264978499
    "J10" = 3847.9913775190598
    "K10" = 2.6821563207204799
    "L10" = 1.65163444500919
    "C11" = 83.584225951514597
    "D11" = 2.1303242954139798
    "E11" = 6.1922344456176397
    "F11" = 0.048329204840320002
    "D12" = 0.85915402787586004
    "F12" = 0.01397050771232
    "C13" = 50.718023469329701
    "D13" = 1.65758338695817
    "E13" = 51.027310954075197
    "F13" = 0.76334992050674999
    "G13" = 232538.11214330501
    "H13" = 7379.8939258723103
    "I13" = 259494.32589182301
    "J13" = 3967.7522524897699
    "K13" = 3.0001551358001501
    "L13" = 1.6410370404201899
    "C15" = 21.856269520457499
    "D15" = 1.6062095914332499
    "E15" = 2.5031687896293202
    "F15" = 0.0043903934187000003
    "C18" = 86.481807284811794
    "D18" = -1.661496877367
    "E18" = 8.8479068632792401
    "G18" = 9932.4355666606298
    "I18" = 1016.18210324762
    "K18" = -1.8932839123123
    "L18" = 0
    "C19" = 55.626879165615797
    "D19" = 0.61335221683148999
    "E19" = 2.2981611901883299
    "C22" = 73.846409214276903
    "E22" = 0.88221876088604001
    "C23" = 43.7035624980184
    "D23" = 0.83032967388839996
    "E23" = 3.02256898899472
    "F23" = 0.00357704455989
    "G23" = 46309.255923861601
    "I23" = 3162.6166979607001
    "K23" = 1.6901736709566999
    "L23" = 0.10596953031502999
    "D24" = 0.17425804987530999
    "C25" = 59.863961322356403
    "D25" = 3.7212473430719299
    "C28" = 11.8097770035551
    "D28" = 0.60995922514342005
    "E28" = 1.4622214478430999
    "F28" = -0.064284723748600003
    "G28" = 11917.482169827499
    "I28" = 1475.5569074474299
    "K28" = 5.4926258838881097
    "L28" = -4.0411496659382999
    "C29" = 31.419536309169001
    "D29" = 2.5545615694441999
    "E29" = 5.4685373611094299
    "C30" = 2902.4325156754398
    "D30" = 43.391111117536397
    "E30" = 2488.7051022872402
    "F30" = 54.276159618058301
    "C31" = 26.610610460810801
    "D31" = 0.40674009918411003
    "E31" = 0.72390877482666005
    "C32" = 286.21617827409
    "D32" = 3.0054320533022199
    "E32" = 133.46488558303199
    "F32" = 0.23815209612037
    "H32" = 41.474962335570602
    "D34" = 5.6184603886397602
    "C35" = 84.747701047379294
    "D35" = 0.078499897059530005
    "F35" = 0.000019259899100000001
    "C36" = 23.818157706091601
    "D36" = 0.83154625612554001
    "E36" = 3.0913713380919501
    "C37" = 40.769348007474903
    "D37" = 0.88780966022451002
    "E37" = 0.45785937248718001
    "F37" = 0.00086494859495000001
    "I37" = 173.652324203213
    "L37" = 0.18909016742358001
    "C38" = 133.32111122870799
    "D38" = 2.1298950766537299
    "E38" = 110.216282415512
    "F38" = 2.26530853347536
    "G38" = 415773.91720674402
    "I38" = 290041.99841353198
    "K38" = 1.70706635322183
    "L38" = 2.1714213826907098
    "C39" = 13.304598115491199
    "D39" = 2.5837811515500602
    "E39" = 1.75096825819927
    "F39" = 0.093769540787730005
    "G39" = 19378.945431099801
    "I39" = 2550.3903261627302
    "K39" = 19.002002995243199
    "L39" = 5.3707470804061597
    "C40" = 34.632141690085597
    "D40" = 1.7895300960721501
    "E40" = 2.2171361348931402
    "F40" = 0.026220508463499999
    "D41" = 0.67688241938577998
    "F41" = 0.11189019841862
    "D42" = 0.96925536724038996
    "C43" = 59.309963804965797
    "D43" = 0.83147767063824995
    "E43" = 6.0347521326413203
    "F43" = 0.34119153773464
    "C44" = 89.222701001299001
    "D44" = 2.0740113798638098
    "E44" = 1.5265352660848099
    "F44" = 0.051321605095570003
    "C45" = 36.354244933545601
    "D45" = 1.78862425205158
    "E45" = 2.5755316338571501
    "F45" = 0.084568786928219994
    "G45" = 261785.89591293101
    "I45" = 19004.291531400799
    "K45" = 4.8970469700439203
    "L45" = 3.66370334676806
    "C46" = 21.8793160002677
    "D46" = 0.31204813150406002
    "E46" = 3.4499683819240499
    "F46" = 0.083614131470300002
    "G46" = 3429.3639898819602
    "K46" = 1.5290284270163601
    "C47" = 21.031136109103301
    "D47" = 2.9499058068541402
    "E47" = 1.13576976930165
    "F47" = 0.12089075828653
    "C48" = 132.00066132218001
    "D48" = 0.88744868624232998
    "E48" = 4.4771787151757296
    "F48" = 0.068592439679869999
    "G48" = 2469.7323733379799
    "I48" = 83.768013760937805
    "K48" = 0.71363146911955
    "L48" = 1.64822738732047
    "C49" = 19.3366716709711
    "D49" = 2.0370532559643002
    "E49" = 3.9707597695248098
    "F49" = 0.20281250637066001
    "G49" = 12196.412289748299
    "I49" = 2504.5170170300798
    "K49" = 11.459731236258101
    "L49" = 5.3969335695398701
    "D50" = 0.78850193477059005
    "C51" = 19.045134926372601
    "D51" = 3.5157631658985702
    "E51" = 0.83367106552458003
    "G51" = 13333.118059254901
    "I51" = 583.63643955244595
    "K51" = 17.676796811222601
    "D52" = 0.17467527858613
    "F52" = 0.00092270565690000001
    "C53" = 22.4954481766427
    "D53" = 0.88789353786990999
    "E53" = 0.26040587938262999
    "F53" = 0.02166281226162
    "C54" = 19.437041885818498
    "D54" = 1.7312049295671601
    "F54" = 0
    "G54" = 590.30296207230799
    "H54" = 52.576693710954501
    "J54" = 0
    "K54" = 14.419764630219699
    "C55" = 28.954427081186299
    "D55" = 4.04414029247447
    "E55" = 0.30887882786839999
    "F55" = 0.02226317852111
    "C56" = 64.561128421569293
    "D56" = 1.0596559472089
    "E56" = 2.23974826202328
    "F56" = 0.039894522587290003
    "G56" = 8886.1937159448007
    "I56" = 308.27895078488399
    "K56" = 1.6781451792839599
    "L56" = 1.85719542058295
    "C57" = 20.262668743872599
    "D57" = 1.13157674592989
    "E57" = 2.8576275607035502
    "F57" = 0.12085091337157
    "G57" = 87012.56
    "I57" = 12271.31
    "K57" = 5.6621477786105396
    "L57" = 4.2384673263219197
    "C58" = 31.944624206945701
    "D58" = 2.1857118210658202
    "E58" = 3.7914454759834499
    "F58" = 0.79392826039983999
    "G58" = 7837.29410293206
    "H58" = 536.24253818028797
    "I58" = 930.19323307777904
    "J58" = 194.782359406496
    "K58" = 7.9562801112408401
    "L58" = 22.219658328377399
    "D59" = 0.57570166368495002
    "C60" = 31.250965434121898
    "D60" = 2.0004036152731501
    "E60" = 45.9751101054581
    "F60" = 1.9375687103094299
    "C61" = 21.994211871909599
    "D61" = 1.5883814484765
    "E61" = 3.1246654323524399
    "F61" = 0.16500143146129001
    "G61" = 147295.77622615601
    "H61" = 10872.7070943935
    "I61" = 20887.055569903201
    "J61" = 1031.02594497268
    "K61" = 7.5481662872294297
    "L61" = 5.1386173238848096
    "C62" = 48.392799259772602
    "D62" = 1.7190782133514899
    "E62" = 23.790392654699101
    "F62" = 0.48201047104340999
    "G62" = 1103702.9574130001
    "H62" = 39227.071848793297
    "I62" = 592590.28810461902
    "J62" = 11852.610139872801
    "K62" = 3.4996233140414201
    "L62" = 2.0921708191330199
    "C63" = 58.391016416880099
    "D63" = 2.2455257375493902
    "E63" = 51.872106364866099
    "F63" = 2.0965902561441201
    "G63" = 53398933.188041098
    "H63" = 1525384.58007487
    "I63" = 51346536.466706403
    "J63" = 1734139.15352162
    "K63" = 3.2812599248645098
    "L63" = 3.7693306021048398
    "C64" = 54.186027277394999
    "D64" = 2.2972676472367
    "E64" = 23.0503446371786
    "F64" = 0.56749137725262999
    "G64" = 2383616.5156442798
    "H64" = 103224.390287897
    "I64" = 997705.09311723104
    "J64" = 23507.5661466452
    "K64" = 4.2457009558126204
    "L64" = 2.4148754221296098
    "C65" = 22.886309786687701
    "D65" = 1.3090991641071099
    "E65" = 14.522904217305801
    "F65" = 0.29792562887295998
    "G65" = 4735341.0003148401
    "H65" = 82423.815100897002
    "I65" = 3263896.45195165
    "K65" = 5.4668051340501798
    "L65" = 2.8936892721296998
    "C66" = 57.897593002498702
    "D66" = 2.2124225516746199
    "E66" = 50.525929732231802
    "F66" = 2.0019414353972902
    "G66" = 54502636.145454101
    "H66" = 1564611.65192367
    "I66" = 51939126.754811101
    "J66" = 1745991.7636614901
    "K66" = 3.2862968708136999
    "L66" = 3.7490072849100402
    "C67" = 74.698667989209
    "D67" = 1.89999089315731
    "E67" = 33.554350961008502
    "F67" = 0.63572462824766995
    "G67" = 615297.22294902301
    "I67" = 299756.62590151699
    "K67" = 2.5379156494075601
    "L67" = 2.1400159900093301
    "C68" = 35.363638346935602
    "D68" = 1.5160661014627901
    "E68" = 2.9548939779057402
    "F68" = 0.11333155094470999
    "G68" = 417760.196137705
    "H68" = 20089.015896853001
    "I68" = 37257.188761380297
    "J68" = 1618.4310084384399
    "K68" = 4.8255281701036701
    "L68" = 4.3732707937223303
    "C69" = 27.2731593759198
    "D69" = 0.74560397091069996
    "E69" = 1.822559835521
    "F69" = -0.041309762541899997
    "G69" = 72793.918420763293
    "I69" = 4868.5808406246797
    "K69" = 2.6196781544193302
    "L69" = -3.4843022781727
    "C70" = 40.009529387521297
    "D70" = 2.1534905252283001
    "E70" = 5.3328770803803902
    "F70" = 0.11254899920245
    "G70" = 70611.408381306304
    "I70" = 8076.3720995988097
    "K70" = 4.39785968082208
    "L70" = 1.6913103497097299
    "C71" = 21.994211871909599
    "D71" = 1.5883814484765
    "E71" = 3.1246654323524399
    "F71" = 0.16500143146129001
    "G71" = 147295.77622615601
    "H71" = 10872.7070943935
    "I71" = 20887.055569903201
    "J71" = 1031.02594497268
    "K71" = 7.5481662872294297
    "L71" = 5.1386173238848096
    "C72" = 42.852043377603401
    "D72" = 0.56623162046838005
    "E72" = 1.1597387061798901
    "F72" = -0.029252814994099999
    "G72" = 58497.5775537503
    "I72" = 1649.20923165064
    "K72" = 1.92898063016232
    "L72" = -3.6285583150638998
    "C73" = 97.150580575442504
    "D73" = 2.2282457091057202
    "E73" = 83.558987780400201
    "F73" = 1.8317643571468001
    "G73" = 614201.53842931602
    "H73" = 13547.723479873701
    "I73" = 548932.83720494597
    "J73" = 10159.0637094711
    "K73" = 2.1640566381752002
    "L73" = 1.94599556927998
    "C74" = 38.914150910211603
    "D74" = 1.7874513915700301
    "E74" = 3.1082864401049002
    "F74" = 0.16197835362891
    "G74" = 129309.295912931
    "I74" = 10523.1915314008
    "K74" = 4.49856880274569
    "L74" = 5.6989449055346597
    "C75" = 83.780145364130604
    "D75" = 3.2706992231559702
    "E75" = 6.9593634424295496
    "F75" = 0.16924942761577999
    "G75" = 71056.0315304658
    "I75" = 4950.7215819058501
    "K75" = 3.4530615736008898
    "L75" = 2.4572915599825502
    "C76" = 80.727632667053101
    "D76" = 1.85463887799892
    "E76" = 50.596975146750097
    "F76" = 0.51174946250520004
    "G76" = 2874960.5508562201
    "I76" = 1864202.1092180701
    "C77" = 53.477096954008204
    "D77" = 2.3635419257264898
    "E77" = 25.852847049598601
    "F77" = 0.65814959416725005
    "G77" = 1540142.5294571701
    "H77" = 70305.5393455718
    "I77" = 730130.252748381
    "J77" = 17439.159955196799
    "K77" = 4.4470380273032202
    "L77" = 2.4463404321988702
    "C78" = 136.117953747744
    "D78" = 4.7886695314059002
    "E78" = 145.48033880186699
    "F78" = 4.4438450054255298
    "G78" = 22614103.480012499
    "H78" = 779935.903648493
    "I78" = 25135268.7084709
    "J78" = 773402.30944887595
    "K78" = 3.5787832368909598
    "L78" = 3.16335678108014
    "C79" = 79.615668389179703
    "D79" = 2.3309947598268299
    "E79" = 79.954929584762496
    "F79" = 2.8403217742048801
    "G79" = 43237368.387037002
    "H79" = 1245871.5255722399
    "I79" = 44829118.4908932
    "J79" = 1628750.5454281799
    "K79" = 3.0134095562729901
    "L79" = 3.7487616886763
    "C80" = 23.336312271310899
    "D80" = 2.50190668486387
    "E80" = 3.9461098479499301
    "F80" = 0.11851745510764999
    "G80" = 53376.168426364799
    "I80" = 8490.4526051462908
    "K80" = 11.016083078642501
    "L80" = 3.0412630763752801
    "C81" = 39.548111042519601
    "D81" = 2.94478126554035
    "E81" = 25.979732466312701
    "F81" = 1.5656339990519099
    "G81" = 1063689.6115458601
    "H81" = 85903.937099901901
    "I81" = 687990.96797804604
    "J81" = 48754.330251442203
    "K81" = 7.3678619457292003
    "L81" = 6.8381018502078597
    "C82" = 51.7289351748109
    "D82" = 1.5966744758071401
    "E82" = 26.535062680877999
    "F82" = 0.54186239235920997
    "G82" = 1050326.78898663
    "H82" = 32747.563619640801
    "I82" = 584099.83549947303
    "J82" = 11552.3060129966
    "K82" = 3.1205424487846698
    "L82" = 2.0782540190715499
    "C83" = 59.592149500822401
    "D83" = 2.1826226310474302
    "E83" = 53.480582971588802
    "F83" = 2.1433543447056098
    "G83" = 52335243.5764952
    "H83" = 1439480.6429749699
    "I83" = 50658545.498728402
    "J83" = 1685384.82327018
    "K83" = 3.1836262936939002
    "L83" = 3.72201777255543
    "C84" = 68.712769420015903
    "D84" = 1.1714342240090001
    "E84" = 3.49794173210659
    "F84" = 0.070335510183060002
    "G84" = 162536.684172387
    "H84" = 3617.5932136899501
    "I84" = 7416.6645466858499
    "K84" = 2.0126602715894801
    "L84" = 3.07582059637953
    "C85" = 7.9069590139872901
    "D85" = 0.036226038133819999
    "E85" = 0.66457435747100002
    "F85" = 0.00002573591555
    "C86" = 29.455500772137299
    "D86" = 1.7225199965432101
    "E86" = 2.9098321107232601
    "F86" = 0.10076536853593999
    "G86" = 460403.59486792999
    "H86" = 26509.355886608799
    "I86" = 46450.138771890299
    "J86" = 1583.89680398928
    "K86" = 5.8442701700256396
    "L86" = 3.6975836626657199
    "C87" = 20.282909856235101
    "D87" = 1.2385336040142401
    "E87" = 6.9743718759155504
    "F87" = 0.19476169230181001
    "G87" = 1074630.9630227201
    "H87" = 66606.029334768798
    "I87" = 371322.20974207798
    "J87" = 10456.453294916901
    "K87" = 6.33568410135946
    "L87" = 2.7820342068002999
    "C88" = 127.38945843560001
    "D88" = 2.1509337962716
    "E88" = 143.39571412978
    "F88" = 2.3926144609662301
    "G88" = 476812.89511249698
    "H88" = 9058.6477861589701
    "I88" = 536881.66936499695
    "J88" = 10070.4764943916
    "K88" = 1.82638582337152
    "L88" = 1.95324293575895
    "C89" = 33.307556015504503
    "D89" = 2.4109216918083698
    "E89" = 19.880443053230401
    "F89" = 0.86659121266951
    "G89" = 6852467.9666192997
    "H89" = 185320.177513203
    "I89" = 4299436.4021371901
    "J89" = 63631.093927947702
    "K89" = 5.0010589642183998
    "L89" = 3.7180435368312001
    "C90" = 87.509539246386197
    "D90" = 2.4919734786051002
    "E90" = 87.694547430089401
    "F90" = 3.0506222933039302
    "G90" = 45474191.896392196
    "H90" = 1273492.5537141999
    "I90" = 46677485.851555601
    "J90" = 1660054.8876155
    "K90" = 3.0545264291933001
    "L90" = 3.77918072473661
    "C91" = 53.998782773172501
    "D91" = 1.58376165638093
    "E91" = 4.4498361265451098
    "F91" = 0.13757758241407
    "G91" = 214398.75226474201
    "H91" = 7237.7598744215102
    "I91" = 15758.2803456291
    "J91" = 512.77818093853898
    "K91" = 3.22751798383598
    "L91" = 4.08106749083721
    "C92" = 13.0854494834518
    "D92" = 0.81918013281608004
    "E92" = 0.80043831162947998
    "F92" = 0.02034778101819
    "H92" = 3953.3138798967698
    "I92" = 3348.1198285922201
    "J92" = 95.069646893458298
    "K92" = 6.3912039735615496
    "L92" = 3.8683507685321099
    "C93" = 2041.6858540900801
    "D93" = 28.382574860605999
    "E93" = 1727.8299513045999
    "F93" = 37.469687938725599
    "G93" = 338019.961329268
    "I93" = 285671.01748942002
    "K93" = 1.5036842087493301
    "L93" = 2.2139923283576399
    "C94" = 377.41637352431502
    "D94" = 3.1295372791010299
    "E94" = 256.63081533232099
    "F94" = 0.060671840758619998
    "G94" = 2087303.4530410799
    "H94" = 5107.2693045811802
    "I94" = 1385341.63284584
    "J94" = 130.06267925722301
    "K94" = 4.8792563887493197
    "L94" = 2.1870440972238798
    "C95" = 42.017340044251704
    "D95" = 1.88166827756084
    "E95" = 2.40866298788068
    "F95" = 0.018155800825949999
    "G95" = 60414.257424916897
    "H95" = 3345.7021388744001
    "I95" = 3419.8121678069701
    "J95" = 28.0694171629005
    "K95" = 4.5520376035502697
    "L95" = 0.72168376043090998
    "C96" = 54.771569875270202
    "D96" = 3.70315425640353
    "E96" = 8.5647284577031702
    "F96" = 0.30919946063092002
    "G96" = 249475.40492598101
    "H96" = 19326.144802024399
    "I96" = 36180.882080076502
    "J96" = 1267.3965109854801
    "K96" = 7.2006550724478897
    "L96" = 1.32447366386932
    "C97" = 30.163745571013099
    "D97" = 1.47153967552399
    "E97" = 3.5068430220568301
    "F97" = 0.10320223518979001
    "G97" = 266709.12833881
    "H97" = 14350.9297191145
    "I97" = 30828.2536908207
    "J97" = 962.51203225982204
    "K97" = 4.2142269043778002
    "L97" = 3.0860156484042398
    "C98" = 15.560854965891499
    "D98" = 1.05020940046308
    "E98" = 1.1703202951081599
    "F98" = 0.027554612200070001
    "G98" = 137431.88391894699
    "H98" = 9715.9932000691497
    "I98" = 11012.4769014965
    "J98" = 382.64407660807899
    "K98" = 6.84990250165162
    "L98" = 3.7483864031178298
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

$ws.Range("A106").Value = "Source : Statistiques de la balance des paiements et de la position extérieure globale du FMI  (mis à jour le 25/10/2023)."

Write-Host "Applied data refresh changes"